# Apply targeted updates to column F ("dSF") for specific rows.
# These correspond to a "repull data, push all data, mean calculation" style
# update where a handful of dSF values were recomputed/corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    4  = 0
    9  = 1
    19 = 0
    34 = 0
    35 = 1
    38 = 3
    40 = -2
    41 = 3
    49 = 2
    52 = -3
    54 = 8
    56 = -4
    58 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
